# Update "想去人数" (want-to-go count) figures in the F column of the
# "展览" and "全部类型" sheets to the newly scraped values.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibition) ----
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value  = 470
$wsExpo.Range("F4").Value  = 7870
$wsExpo.Range("F6").Value  = 216
$wsExpo.Range("F13").Value = 446
$wsExpo.Range("F15").Value = 70
$wsExpo.Range("F17").Value = 5764
$wsExpo.Range("F18").Value = 169
$wsExpo.Range("F19").Value = 243
$wsExpo.Range("F20").Value = 1507
$wsExpo.Range("F22").Value = 356

# ---- Sheet "全部类型" (All types) ----
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 470
$wsAll.Range("F4").Value  = 7870
$wsAll.Range("F6").Value  = 216
$wsAll.Range("F13").Value = 446
$wsAll.Range("F15").Value = 70
$wsAll.Range("F18").Value = 5764
$wsAll.Range("F20").Value = 169
$wsAll.Range("F21").Value = 243
$wsAll.Range("F22").Value = 1508
$wsAll.Range("F24").Value = 356
